# Updates cryptos list price / 1h-volume figures to match the latest
# scrape, and fixes the HuobiToken/ImmutableX row ordering (rows 34-35).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.925.22"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "1.862.75"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'305.02"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.5066"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.3625"
$ws.Range("E8").Value = "  -3.38%  "
$ws.Range("D9").Value = "'0.07175"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").Value = "'0.8961"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "1.855.57"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("D14").Value = "'92.57"
$ws.Range("E14").Value = "  +3.46%  "
$ws.Range("D15").Value = "'5.240"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "'0.000008493"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "'0.9997"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "26.959.83"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("D21").Value = "'5.028"
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("D22").Value = "2.089.45"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").Value = "'10.33"
$ws.Range("E23").Value = "  -2.74%  "
$ws.Range("D24").Value = "'6.432"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("D25").Value = "'148.02"
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("D28").Value = "'2.062"
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("D29").Value = "'113.35"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").Value = "'4.673"
$ws.Range("E30").Value = "  -2.02%  "
$ws.Range("D31").Value = "'4.677"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "'0.09259"
$ws.Range("E32").Value = "  +2.90%  "
$ws.Range("D33").Value = "'0.05084"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7487"
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.988"
$ws.Range("E35").Value = "  -3.55%  "
$ws.Range("D36").Value = "'1.152"
$ws.Range("E36").Value = "  -0.87%  "
$ws.Range("D37").Value = "'3.276"
$ws.Range("E37").Value = "  +7.66%  "
$ws.Range("D38").Value = "'2.525"
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").Value = "'0.02003"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").Value = "'1.078"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").Value = "'0.5472"
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("D42").Value = "'117.92"
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("D43").Value = "'6.493"
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").Value = "'8.558"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").Value = "'0.1470"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").Value = "'0.4681"
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("D47").Value = "'0.9993"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").Value = "'10.11"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").Value = "'1.565"
$ws.Range("D50").Value = "'36.93"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").Value = "'63.02"
$ws.Range("E51").Value = "  -2.66%  "
